# Actualización desde MV -datos-
# Adds a new column BH with header "Agosto.2021" (a copy of the last
# existing data column, BG "Mayo.2021") to reflect the newest data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column BH (60), mirroring existing header formatting
$ws.Range("BG1").Copy()
$ws.Range("BH1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Cells.Item(1, 60).Value = "Agosto.2021"

# Copy last column (BG, 59) values into new column (BH, 60) for each data row
for ($row = 2; $row -le 19; $row++) {
    $lastVal = $ws.Cells.Item($row, 59).Value2
    if ($null -ne $lastVal) {
        $ws.Cells.Item($row, 60).Value = $lastVal
    }
}
